$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 452-453. This shifts the existing rows 452-488
# (Fuerte/Hass/Edranol/Negra de La Cruz entries) down to 454-490, matching
# the new sheet dimension A1:T490.
$ws.Rows.Item(452).Resize(2).Insert()

# New row 452: Comercializadora del Agro de Limarí, Palta, Negra de La Cruz, Primera
$ws.Cells.Item(452, 1).Value = 2
$ws.Cells.Item(452, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(452, 3).Value = "Coquimbo"
$ws.Cells.Item(452, 4).Value = 44783
$ws.Cells.Item(452, 5).Value = 4
$ws.Cells.Item(452, 6).Value = "Fruta"
$ws.Cells.Item(452, 7).Value = 100106
$ws.Cells.Item(452, 8).Value = "Oleaginosos"
$ws.Cells.Item(452, 9).Value = 100106002
$ws.Cells.Item(452, 10).Value = "Palta"
$ws.Cells.Item(452, 11).Value = "Negra de La Cruz"
$ws.Cells.Item(452, 12).Value = "Primera"
$ws.Cells.Item(452, 13).Value = 300
$ws.Cells.Item(452, 14).Value = 1500
$ws.Cells.Item(452, 15).Value = 1600
$ws.Cells.Item(452, 16).Value = 1550
$ws.Cells.Item(452, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(452, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(452, 19).Value = 1550
$ws.Cells.Item(452, 20).Value = 1

# New row 453: Comercializadora del Agro de Limarí, Palta, Negra de La Cruz, Segunda
$ws.Cells.Item(453, 1).Value = 2
$ws.Cells.Item(453, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(453, 3).Value = "Coquimbo"
$ws.Cells.Item(453, 4).Value = 44783
$ws.Cells.Item(453, 5).Value = 4
$ws.Cells.Item(453, 6).Value = "Fruta"
$ws.Cells.Item(453, 7).Value = 100106
$ws.Cells.Item(453, 8).Value = "Oleaginosos"
$ws.Cells.Item(453, 9).Value = 100106002
$ws.Cells.Item(453, 10).Value = "Palta"
$ws.Cells.Item(453, 11).Value = "Negra de La Cruz"
$ws.Cells.Item(453, 12).Value = "Segunda"
$ws.Cells.Item(453, 13).Value = 360
$ws.Cells.Item(453, 14).Value = 1300
$ws.Cells.Item(453, 15).Value = 1400
$ws.Cells.Item(453, 16).Value = 1350
$ws.Cells.Item(453, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(453, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(453, 19).Value = 1350
$ws.Cells.Item(453, 20).Value = 1
